$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly data between row 2 and row 3 for the columns that
# differ: D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), P (Precio $/Kg).
$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $v2 = $cell2.Value2
    $v3 = $cell3.Value2
    $cell2.Value2 = $v3
    $cell3.Value2 = $v2
}
